# Apply the update described by the diff: insert "line7" and "line8" entries
# right after the existing "line6" row, which shifts the "extr1".."extr8"
# rows down by two rows, and update the numeric/boolean data accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 8-17 (columns A-E). Rows 2-7 (line1..line6)
# are unchanged by this edit.
$data = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $false },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $false },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $false },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $false },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $true  }
)

# Row 7 ("line6") already carries the bold / centered / bordered style used
# by every cell in column A (A2:A15). Copy that formatting onto the two
# brand-new rows (16 and 17) so their "A" cell matches the rest of the table.
$ws.Cells.Item(7, 1).Copy() | Out-Null

foreach ($item in $data) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.A
    if ($r -gt 15) {
        $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    }

    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}

$excel.CutCopyMode = 0
